function Set-ParaText($para, $text) {
    # Assigning to Range.Text only overwrites the first run when a
    # paragraph holds several runs (e.g. ones split by proofErr spans),
    # so rebuild a range that excludes the trailing paragraph mark and
    # retext that instead - this clears every run in one shot.
    $doc = $para.Range.Document
    $s = $para.Range.Start
    $e = $para.Range.End - 1
    $rng = $doc.Range($s, $e)
    $rng.Text = $text
}

$d = $word.ActiveDocument

# ------------------------------------------------------------------
# 1) Move the "Créer un dossier Upload ..." list paragraph (currently
#    paragraph 6) to sit right after "Récupérer les fichiers"
#    (paragraph 3), then retext it to "Chiffrer".
# ------------------------------------------------------------------
$creerPara = $d.Paragraphs.Item(6)
$creerPara.Range.Cut()

$afterRecuperer = $d.Paragraphs.Item(3).Range.End
$d.Range($afterRecuperer, $afterRecuperer).Paste()

Set-ParaText $d.Paragraphs.Item(4) "Chiffrer"

# ------------------------------------------------------------------
# 2) The paragraph that used to hold "Chaque dossier utilisateur
#    nommer par « nom_prénom_id » contiendra les fichiers uploaders "
#    is now paragraph 7. Split its role in two:
#      - paragraph 7 becomes the (merged, single-run) "Créer un
#        dossier Upload qui contiendra Tout les dossiers
#        utilisateurs" text.
#      - a new paragraph 8 (same list style) takes the merged
#        "Chaque dossier utilisateur..." text.
# ------------------------------------------------------------------
$chaquePara = $d.Paragraphs.Item(7)
$chaquePara.Range.InsertParagraphAfter()

Set-ParaText $d.Paragraphs.Item(8) "Chaque dossier utilisateur nommer par « nom_prénom_id » contiendra les fichiers uploaders "
Set-ParaText $d.Paragraphs.Item(7) "Créer un dossier Upload qui contiendra Tout les dossiers utilisateurs"
